# Refresh the cryptos list (Price in column D, Volume(1h) in column E)
# to match the latest scrape, per the commit "Updated cryptos list on
# Thu Jul 11 22:37:31 UTC 2024 with GitHub Actions".
#
# Column D values are plain text in the source sheet (things like
# "57.505.41" or "1.00" are literal strings, not numbers - some even use
# a thousands "." so they cannot be numbers at all). A direct
# `Range.Value = "526.19"` assignment would get auto-coerced to a Number
# by Excel, which would corrupt the text. To keep the new price a literal
# string we round-trip it through a text-producing formula (`="526.19"`)
# and then flatten the formula to its plain value via Copy/PasteSpecial
# (values-only), which keeps the cell a text cell instead of a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="57.505.41"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").Formula = '="3.107.69"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = '  +0.14%  '

$ws.Range("D5").Formula = '="526.19"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  +0.56%  '

$ws.Range("D6").Formula = '="137.01"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  -3.34%  '

$ws.Range("D7").Formula = '="0.999"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").Formula = '="3.103.92"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").Formula = '="0.447"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  +2.03%  '

$ws.Range("D10").Formula = '="7.30"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  +1.22%  '

$ws.Range("E11").Value = '  -0.66%  '

$ws.Range("D12").Formula = '="0.395"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  +2.68%  '

$ws.Range("D13").Formula = '="3.641.78"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = '  +0.13%  '

$ws.Range("E14").Value = '  +2.94%  '

$ws.Range("D15").Formula = '="25.33"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  -1.64%  '

$ws.Range("E16").Value = '  -0.54%  '

$ws.Range("D17").Formula = '="57.587.05"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = '  -0.42%  '

$ws.Range("D18").Formula = '="3.099.78"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  -0.09%  '

$ws.Range("D19").Formula = '="5.94"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  -2.55%  '

$ws.Range("D20").Formula = '="12.37"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  -3.14%  '

$ws.Range("D21").Formula = '="7.87"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  -2.12%  '

$ws.Range("D22").Formula = '="345.36"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  +2.22%  '

$ws.Range("E23").Value = '  -0.12%  '

$ws.Range("D24").Formula = '="67.69"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  +1.70%  '

$ws.Range("E25").Value = '  -2.12%  '

$ws.Range("D26").Formula = '="0.167"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  -1.61%  '

$ws.Range("D27").Formula = '="0.999"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  -0.27%  '

$ws.Range("D28").Formula = '="0.0₃0893"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  -2.30%  '

$ws.Range("D29").Formula = '="7.43"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  +3.63%  '

$ws.Range("E30").Value = '  +0.00%  '

$ws.Range("E31").Value = '  +0.29%  '

$ws.Range("D32").Formula = '="6.03"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  -6.94%  '

$ws.Range("D33").Formula = '="20.80"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  -0.49%  '

$ws.Range("D34").Formula = '="4.96"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  +7.69%  '

$ws.Range("E35").Value = '  -3.23%  '

$ws.Range("D36").Formula = '="158.40"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  +1.56%  '

$ws.Range("D37").Formula = '="6.06"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  -0.79%  '

$ws.Range("D38").Formula = '="25.89"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  -4.10%  '

$ws.Range("E39").Value = '  -1.58%  '

$ws.Range("D40").Formula = '="1.61"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  +5.73%  '

$ws.Range("D41").Formula = '="0.0660"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  -0.01%  '

$ws.Range("E42").Value = '  +3.73%  '

$ws.Range("E43").Value = '  +2.27%  '

$ws.Range("D44").Formula = '="3.147.47"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +0.14%  '

$ws.Range("D45").Formula = '="2.379.72"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  +3.51%  '

$ws.Range("D46").Formula = '="36.61"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  -0.50%  '

$ws.Range("D47").Formula = '="1.00"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  -0.01%  '

$ws.Range("D48").Formula = '="0.0267"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  +3.24%  '

$ws.Range("D49").Formula = '="0.974"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  -0.85%  '

$ws.Range("D50").Formula = '="5.98"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  -0.66%  '

$ws.Range("D51").Formula = '="19.79"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  -3.47%  '
